$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.231.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.021"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.56%  "
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.019"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4789"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3721"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07322"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9353"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.35"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07866"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.536"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008732"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.267.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.104"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.85"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.000"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.988"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08886"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.350"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.187"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.587"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7404"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.124"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02032"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.94%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.997"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05262"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5329"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.110"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1527"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.330"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.61"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4789"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.021"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.635"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "66.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06076"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.27%  "
